$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 0.1972685887708649
$ws.Range("C2").Value = 0.5508345978755691
$ws.Range("J2").Value = 0.01213960546282246
$ws.Range("P2").Value = 0.1684370257966616
$ws.Range("S2").Value = 0.07132018209408195

# Row 3
$ws.Range("B3").Value = 0.007556675062972292
$ws.Range("C3").Value = 0.03778337531486146
$ws.Range("J3").Value = 0.0327455919395466
$ws.Range("P3").Value = 0.743073047858942
$ws.Range("S3").Value = 0.1788413098236776

# Row 4
$ws.Range("J4").Value = 0.08571428571428572
$ws.Range("O4").Value = 0.01904761904761905
$ws.Range("P4").Value = 0.5047619047619047
$ws.Range("S4").Value = 0.3904761904761905

# Row 6
$ws.Range("B6").Value = 0.06776180698151951
$ws.Range("D6").Value = 0.01232032854209446
$ws.Range("F6").Value = 0.05544147843942505
$ws.Range("J6").Value = 0.2197125256673511
$ws.Range("O6").Value = 0.02258726899383984
$ws.Range("Q6").Value = 0.1806981519507187
$ws.Range("R6").Value = 0.06365503080082136
$ws.Range("S6").Value = 0.37782340862423

# Row 7
$ws.Range("B7").Value = 0.1226666666666667
$ws.Range("D7").Value = 0.02133333333333333
$ws.Range("E7").Value = 0.002666666666666667
$ws.Range("F7").Value = 0.06133333333333333
$ws.Range("J7").Value = 0.1226666666666667
$ws.Range("O7").Value = 0.01333333333333333
$ws.Range("Q7").Value = 0.184
$ws.Range("R7").Value = 0.05066666666666667
$ws.Range("S7").Value = 0.4213333333333333

# Row 8
$ws.Range("B8").Value = 0.07916666666666666
$ws.Range("D8").Value = 0.01875
$ws.Range("E8").Value = 0.001041666666666667
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.10625
$ws.Range("O8").Value = 0.01666666666666667
$ws.Range("Q8").Value = 0.209375
$ws.Range("R8").Value = 0.08020833333333334
$ws.Range("S8").Value = 0.421875

# Row 9
$ws.Range("B9").Value = 0.1142857142857143
$ws.Range("D9").Value = 0.03516483516483516
$ws.Range("F9").Value = 0.06813186813186813
$ws.Range("J9").Value = 0.08571428571428572
$ws.Range("O9").Value = 0.008791208791208791
$ws.Range("Q9").Value = 0.2131868131868132
$ws.Range("R9").Value = 0.07912087912087912
$ws.Range("S9").Value = 0.3956043956043956

# Row 10
$ws.Range("B10").Value = 0.1239636794315041
$ws.Range("D10").Value = 0.02210817212791157
$ws.Range("E10").Value = 0.0007895775759968417
$ws.Range("F10").Value = 0.07816818002368733
$ws.Range("J10").Value = 0.1014607185155942
$ws.Range("O10").Value = 0.02092380576391631
$ws.Range("Q10").Value = 0.2143703118831425
$ws.Range("R10").Value = 0.07698381365969206
$ws.Range("S10").Value = 0.3612317410185551

# Row 11
$ws.Range("G11").Value = 0.09411764705882353
$ws.Range("J11").Value = 0.1176470588235294
$ws.Range("K11").Value = 0.1313725490196078
$ws.Range("L11").Value = 0.6431372549019608
$ws.Range("S11").Value = 0.01372549019607843

# Row 12
$ws.Range("G12").Value = 0.7838616714697406
$ws.Range("J12").Value = 0.1440922190201729
$ws.Range("K12").Value = 0.005763688760806916
$ws.Range("L12").Value = 0.04610951008645533
$ws.Range("S12").Value = 0.02017291066282421

# Row 13
$ws.Range("G13").Value = 0.62
$ws.Range("J13").Value = 0.32
$ws.Range("S13").Value = 0.06

# Row 14
$ws.Range("G14").Value = 0.8333333333333334
$ws.Range("J14").Value = 0.1666666666666667

# Row 15
$ws.Range("F15").Value = 0.01204819277108434
$ws.Range("H15").Value = 0.1827309236947791
$ws.Range("I15").Value = 0.06425702811244979
$ws.Range("J15").Value = 0.3775100401606425
$ws.Range("K15").Value = 0.04819277108433735
$ws.Range("M15").Value = 0.01204819277108434
$ws.Range("O15").Value = 0.07228915662650602
$ws.Range("S15").Value = 0.2309236947791165

# Row 16
$ws.Range("F16").Value = 0.02314814814814815
$ws.Range("H16").Value = 0.1574074074074074
$ws.Range("I16").Value = 0.09027777777777778
$ws.Range("J16").Value = 0.3981481481481481
$ws.Range("K16").Value = 0.1134259259259259
$ws.Range("M16").Value = 0.03472222222222222
$ws.Range("O16").Value = 0.04398148148148148
$ws.Range("S16").Value = 0.1388888888888889

# Row 17
$ws.Range("F17").Value = 0.01780415430267062
$ws.Range("H17").Value = 0.181998021760633
$ws.Range("I17").Value = 0.1107814045499505
$ws.Range("J17").Value = 0.4164193867457963
$ws.Range("K17").Value = 0.0771513353115727
$ws.Range("M17").Value = 0.02176063303659743
$ws.Range("N17").Value = 0.0009891196834817012
$ws.Range("O17").Value = 0.08209693372898121
$ws.Range("S17").Value = 0.09099901088031652

# Row 18
$ws.Range("F18").Value = 0.0223463687150838
$ws.Range("H18").Value = 0.1787709497206704
$ws.Range("I18").Value = 0.1201117318435754
$ws.Range("J18").Value = 0.4189944134078212
$ws.Range("K18").Value = 0.0893854748603352
$ws.Range("M18").Value = 0.01396648044692737
$ws.Range("N18").Value = 0.002793296089385475
$ws.Range("O18").Value = 0.08379888268156424
$ws.Range("S18").Value = 0.06983240223463687

# Row 19
$ws.Range("F19").Value = 0.01521826191429716
$ws.Range("H19").Value = 0.2226672006407689
$ws.Range("I19").Value = 0.09291149379255106
$ws.Range("J19").Value = 0.3704445334401282
$ws.Range("K19").Value = 0.1029235082098518
$ws.Range("M19").Value = 0.0224269122947537
$ws.Range("N19").Value = 0.001601922306768122
$ws.Range("O19").Value = 0.07889467360833
$ws.Range("S19").Value = 0.09291149379255106
